$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header: A1 changes from numeric 0 to text "Description"
$ws.Range("A1").Value = "Description"

# Row 2
$ws.Range("A2").Value = "Move Robot2 to location (2, 8) and remove the toolkit."
$ws.Range("E2").Value = $false

# Row 3
$ws.Range("A3").Value = "Move Robot26 to location (11, 4) and remove the liquid spill."

# Row 4
$ws.Range("A4").Value = "Move Robot42 to location (9, 5) and remove the large debris."
$ws.Range("B4").Value = $true

# Row 5
$ws.Range("A5").Value = "Move Robot48 to location (5, 6) and remove the dust."
$ws.Range("C5").Value = $false

# Row 6
$ws.Range("A6").Value = "Move Robot31 to location (9, 4) and remove the grass."

# Row 7
$ws.Range("A7").Value = "Move Robot8 to location (8, 12) and remove the small debris."
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = $false

# Row 8
$ws.Range("A8").Value = "Move Robot23 to location (11, 1) and remove the vehicle."

# Row 9
$ws.Range("A9").Value = "Move Robot23 to location (12, 10) and remove the construction materials."

# Row 10
$ws.Range("A10").Value = "Move Robot14 to location (7, 11) and remove the tree branches."

# Row 11
$ws.Range("A11").Value = "Move Robot15 to location (5, 3) and remove the screws."
